# Add "Traversal Time" column (T) to Sheet1, computing each block's
# traversal time from the already-derived speed (MPH, column R) and
# block length in yards (column S).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header
$ws.Range("T1").Value = "Traversal Time"

# Formula for every data row (2 through 78). Mirror the existing R/S
# column shared-formula grouping (2 / 3:66 / 67:78) so the new column
# follows the same structure as its neighbours.
$ws.Range("T2").Formula = "= S2/(R2*0.488889)"
$ws.Range("T3:T66").Formula = "= S3/(R3*0.488889)"
$ws.Range("T67:T78").Formula = "= S67/(R67*0.488889)"

# Give the new column a sensible display width (matches the author's
# auto-fit result of 13 character-widths).
$ws.Columns.Item(20).ColumnWidth = 12.166666666666666

# Restore the view: scrolled so column H is first visible, with the new
# cell V69 selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 8
$ws.Range("V69").Select()
